# SWD May-2023 challenge: add one more bullet to the "Internal Discussion"
# list and nudge the window/selection to reflect having scrolled down to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New discussion bullet appended below "Compare simple model to ODE" (B12).
$ws.Range("B13").Value = "This would be a straight line versus a curve."

# Reflect the resulting view state: window scrolled further down and the
# selection moved to the next empty row beneath the new bullet.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$win.Width = 11844
$win.Height = 7944
$win.Top = 60
$win.Left = 24

[void]$ws.Range("B14").Select()
